$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the sheet dimension-relevant content: clear stale cells whose data moved elsewhere
$ws.Range("H5,H6,H7,B23,B24,H26,C27,C28,D29,D30,E31,E32,F33,F34,G35,G36,B37,B38,H39,C41,C42,D43,D44,E45,E46,F47,F48,G49,G50").ClearContents()

# 2) Write the new/updated label cells (column A row labels + header row text)
$ws.Range("G1").Value = "d=6"
$ws.Range("H1").Value = "d=7"
$ws.Range("I1").Value = "d=10"
$ws.Range("A20").Value = "ARMA_I(0,6,0)"
$ws.Range("A21").Value = "ARMA_I(0,6,1)"
$ws.Range("A22").Value = "ARMA_I(0,6,2)"
$ws.Range("A23").Value = "ARMA_I(0,7,0)"
$ws.Range("A24").Value = "ARMA_I(0,7,1)"
$ws.Range("A25").Value = "ARMA_I(0,7,2)"
$ws.Range("A26").Value = "ARMA_I(1,1,0)"
$ws.Range("A27").Value = "ARMA_I(1,1,1)"
$ws.Range("A28").Value = "ARMA_I(1,10,0)"
$ws.Range("A29").Value = "ARMA_I(1,10,1)"
$ws.Range("A30").Value = "ARMA_I(1,2,0)"
$ws.Range("A31").Value = "ARMA_I(1,2,1)"
$ws.Range("A32").Value = "ARMA_I(1,3,0)"
$ws.Range("A33").Value = "ARMA_I(1,3,1)"
$ws.Range("A34").Value = "ARMA_I(1,4,0)"
$ws.Range("A35").Value = "ARMA_I(1,4,1)"
$ws.Range("A36").Value = "ARMA_I(1,5,0)"
$ws.Range("A37").Value = "ARMA_I(1,5,1)"
$ws.Range("A38").Value = "ARMA_I(1,6,0)"
$ws.Range("A39").Value = "ARMA_I(1,6,1)"
$ws.Range("A40").Value = "ARMA_I(1,7,0)"
$ws.Range("A41").Value = "ARMA_I(1,7,1)"
$ws.Range("A42").Value = "ARMA_I(2,1,0)"
$ws.Range("A43").Value = "ARMA_I(2,1,2)"
$ws.Range("A44").Value = "ARMA_I(2,10,0)"
$ws.Range("A45").Value = "ARMA_I(2,10,2)"
$ws.Range("A46").Value = "ARMA_I(2,2,0)"
$ws.Range("A47").Value = "ARMA_I(2,2,2)"
$ws.Range("A48").Value = "ARMA_I(2,3,0)"
$ws.Range("A49").Value = "ARMA_I(2,3,2)"
$ws.Range("A50").Value = "ARMA_I(2,4,0)"
$ws.Range("A51").Value = "ARMA_I(2,4,2)"
$ws.Range("A52").Value = "ARMA_I(2,5,0)"
$ws.Range("A53").Value = "ARMA_I(2,5,2)"
$ws.Range("A54").Value = "ARMA_I(2,6,0)"
$ws.Range("A55").Value = "ARMA_I(2,6,2)"
$ws.Range("A56").Value = "ARMA_I(2,7,0)"
$ws.Range("A57").Value = "ARMA_I(2,7,2)"

# 3) Write the new/updated numeric data cells
$ws.Range("I5").Value = 97.71058970844493
$ws.Range("I6").Value = 93.54218735965073
$ws.Range("I7").Value = 96.26978552796047
$ws.Range("G20").Value = 98.81619798517222
$ws.Range("G21").Value = 99.11933192938405
$ws.Range("G22").Value = 97.82771226540473
$ws.Range("H23").Value = 97.33718409582383
$ws.Range("H24").Value = 96.37610055138501
$ws.Range("H25").Value = 98.09989924092648
$ws.Range("B26").Value = -0.6190711210521657
$ws.Range("B27").Value = 0.6841272603041058
$ws.Range("I28").Value = 97.5104499833141
$ws.Range("I29").Value = 96.5701498421321
$ws.Range("C30").Value = -0.2508026900248366
$ws.Range("C31").Value = -0.1531654487793082
$ws.Range("D32").Value = 0.1298021208374179
$ws.Range("D33").Value = -0.2688808903732567
$ws.Range("E34").Value = -0.3808122002884384
$ws.Range("E35").Value = 0.8061354655834801
$ws.Range("F36").Value = 4.260618407249224
$ws.Range("F37").Value = 28.13318486538547
$ws.Range("G38").Value = 98.80431091244064
$ws.Range("G39").Value = 98.43771467268675
$ws.Range("H40").Value = 98.44083755300801
$ws.Range("H41").Value = 98.61522676209977
$ws.Range("B42").Value = -1.242244261924974
$ws.Range("B43").Value = -0.6600914489611678
$ws.Range("I44").Value = 96.08358411283866
$ws.Range("I45").Value = 90.99781590956033
$ws.Range("C46").Value = 0.7658504971464016
$ws.Range("C47").Value = 0.0391701932455986
$ws.Range("D48").Value = -0.04861296320771899
$ws.Range("D49").Value = 0.0497594565789778
$ws.Range("E50").Value = -0.08599285014618893
$ws.Range("E51").Value = 0.1243771677100176
$ws.Range("F52").Value = -0.7516248542411386
$ws.Range("F53").Value = 29.48661215871149
$ws.Range("G54").Value = 96.41618907505205
$ws.Range("G55").Value = 95.89417146762578
$ws.Range("H56").Value = 98.10516031525961
$ws.Range("H57").Value = 98.06278609418102

# 4) Apply the existing bold/centered/bordered header-style to brand-new cells
#    (copy format from a cell that already carries that style, then paste-format only)
$ws.Range("A2").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E34").PasteSpecial(-4122)
$ws.Range("E35").PasteSpecial(-4122)
$ws.Range("F36").PasteSpecial(-4122)
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("G38").PasteSpecial(-4122)
$ws.Range("G39").PasteSpecial(-4122)
$ws.Range("H41").PasteSpecial(-4122)
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B43").PasteSpecial(-4122)
$ws.Range("I44").PasteSpecial(-4122)
$ws.Range("I45").PasteSpecial(-4122)
$ws.Range("C46").PasteSpecial(-4122)
$ws.Range("C47").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E50").PasteSpecial(-4122)
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("E51").PasteSpecial(-4122)
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("G54").PasteSpecial(-4122)
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("G55").PasteSpecial(-4122)
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("H56").PasteSpecial(-4122)
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("H57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "edit complete"
